$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New article data for rows 2-13 (date columns A & D are date-like strings,
# so force text format on those columns first to avoid Excel's automatic
# date coercion when the value is assigned).
$ws.Range("A2:A13").NumberFormat = "@"
$ws.Range("D2:D13").NumberFormat = "@"

$data = @(
    @("2026-01-31","기업","엔비디아 오픈AI 145조 투자 계획 교착 상태","2026-01-31","https://n.news.naver.com/mnews/article/003/0013741164?sid=104"),
    @("2026-01-31","기업","엔비디아, 오픈AI 145조 투자 보류…`"내부 회의론 커져`"","2026-01-31","https://n.news.naver.com/mnews/article/015/0005244994?sid=104"),
    @("2026-01-31","기술","정부, 독자AI 모델 기반 '국방 AI' 개발한다","2026-01-30","https://n.news.naver.com/mnews/article/003/0013740053?sid=105"),
    @("2026-01-31","기술","KT LLM '믿:음' 개발 이끈 최고AI책임자, NC AI로","2026-01-30","https://n.news.naver.com/mnews/article/001/0015876583?sid=105"),
    @("2026-01-31","산업","오세훈 서울시장, 피지컬 AI 거점으로 양재·수서·용산 꼽았다","2026-01-30","https://n.news.naver.com/mnews/article/277/0005714412?sid=102"),
    @("2026-01-31","산업","서울시민 2000명 앞에서 오세훈 시장이 밝힌 피지컬AI 청사진","2026-01-31","https://n.news.naver.com/mnews/article/277/0005714785?sid=102"),
    @("2026-01-31","정책","대웅제약, AI 에이전트 활용해 전공의 수련 환경 지원…`"행정 감소`"","2026-01-30","https://n.news.naver.com/mnews/article/421/0008744150?sid=103"),
    @("2026-01-31","정부(과기부)","아크릴, 실시간 로봇 제어 기술로 피지컬 AI 시장 공략","2026-01-31","http://www.technoa.co.kr/news/articleView.html?idxno=100468"),
    @("2026-01-31","정부(과기부)","기후 대응의 해법을 국제 무대에 올린 한밭대","2026-01-31","https://www.ccdailynews.com/news/articleView.html?idxno=2395017"),
    @("2026-01-31","정부(과기부)","광주시·국토부, 'AI 모빌리티 국가시범도시'본격 추진","2026-01-31","http://www.siminsori.com/news/articleView.html?idxno=310069"),
    @("2026-01-31","정부(과기부)","韓 '국대 AI' 제친 4년차 中 '스타트업' 문샷...돌파구는? [김성태의...","2026-01-31","https://n.news.naver.com/mnews/article/011/0004585136?sid=105"),
    @("2026-01-31","정부(과기부)","경남, SW미래채움 연차평가 최고등급. . . 7년간 171억 투입","2026-01-31","https://www.straightnews.co.kr/news/articleView.html?idxno=293779")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# The old sheet had 13 data rows (rows 2-14); the refreshed list only has 12
# (rows 2-13), so remove the now-unused last row entirely.
$ws.Rows.Item(14).Delete()

Write-Output "updated $($r - 2) rows and removed row 14"
